# Refresh crypto market snapshot: prices + 1h volume deltas.
# Cells in column D/E are stored as literal text in the source data
# (not numbers), so any value that looks numeric is written with the
# cell pre-formatted as Text ('@') to stop Excel's automatic
# "convert numeric-looking text to a Number" behaviour from silently
# dropping significant digits (e.g. trailing zeros). The style is put
# back to Normal immediately afterwards so no visible formatting or
# style-index change is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.214.40'
$ws.Range('E2').Value = '  -0.80%  '
$ws.Range('D3').Value = '3.157.06'
$ws.Range('E3').Value = '  -0.64%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '613.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.63'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.09%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '3.149.93'
$ws.Range('E8').Value = '  -0.84%  '
$ws.Range('E9').Value = '  -0.74%  '
$ws.Range('E10').Value = '  -1.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.43'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.81%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.473'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.80%  '
$ws.Range('E13').Value = '  -0.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.50'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.11%  '
$ws.Range('D15').Value = '3.670.35'
$ws.Range('E15').Value = '  -0.62%  '
$ws.Range('E16').Value = '  +2.79%  '
$ws.Range('D17').Value = '64.174.51'
$ws.Range('E17').Value = '  -0.80%  '
$ws.Range('D18').Value = '3.157.62'
$ws.Range('E18').Value = '  -0.53%  '
$ws.Range('E19').Value = '  -2.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '477.84'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.66%  '
$ws.Range('E21').Value = '  -0.89%  '
$ws.Range('B22').Value = 'Polygon'
$ws.Range('C22').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.714'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.71%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.69'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.48%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.66'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.84'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.01%  '
$ws.Range('E28').Value = '  -1.61%  '
$ws.Range('E29').Value = '  +0.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.118'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.71%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.10'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.30%  '
$ws.Range('E32').Value = '  +0.41%  '
$ws.Range('E33').Value = '  -1.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.31'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.24%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.14'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.97%  '
$ws.Range('D36').Value = '0.0₃0794'
$ws.Range('E36').Value = '  +7.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '53.10'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.83%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.17'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '462.52'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0400'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.120'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.36'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.70%  '
$ws.Range('D44').Value = '2.856.68'
$ws.Range('E44').Value = '  -0.42%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.31'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.93%  '
$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.268'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.44'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '26.52'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.35%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.998'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('E50').Value = '  -1.71%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '119.65'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.50%  '
